# Refresh cryptocurrency price/volume snapshot (GitHub Actions data update).
# Price (column D) and Volume(1h) (column E) values are stored as plain text
# (e.g. "60.906.13", "  -3.56%  "), so a leading apostrophe is used when
# assigning .Value to force Excel to keep them as text instead of silently
# coercing number-looking strings (like "7.97") into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.906.13"
$ws.Range("E2").Value = "'  -3.56%  "
$ws.Range("D3").Value = "'3.354.74"
$ws.Range("E3").Value = "'  -2.93%  "
$ws.Range("E4").Value = "'  +0.01%  "
$ws.Range("D5").Value = "'567.73"
$ws.Range("E5").Value = "'  -1.90%  "
$ws.Range("D6").Value = "'149.80"
$ws.Range("E6").Value = "'  +0.34%  "
$ws.Range("E7").Value = "'  +0.14%  "
$ws.Range("D9").Value = "'7.97"
$ws.Range("E9").Value = "'  +1.64%  "
$ws.Range("E10").Value = "'  -1.46%  "
$ws.Range("D11").Value = "'0.415"
$ws.Range("E11").Value = "'  +1.80%  "
$ws.Range("D12").Value = "'3.928.12"
$ws.Range("E12").Value = "'  -2.96%  "
$ws.Range("E13").Value = "'  +1.22%  "
$ws.Range("D14").Value = "'28.07"
$ws.Range("E14").Value = "'  -1.96%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.0000170"
$ws.Range("E15").Value = "'  -1.32%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "'3.357.23"
$ws.Range("E16").Value = "'  -2.92%  "
$ws.Range("D17").Value = "'60.986.56"
$ws.Range("E17").Value = "'  -3.38%  "
$ws.Range("D18").Value = "'6.32"
$ws.Range("E18").Value = "'  -2.30%  "
$ws.Range("D19").Value = "'14.17"
$ws.Range("E19").Value = "'  -2.07%  "
$ws.Range("E20").Value = "'  -3.44%  "
$ws.Range("D21").Value = "'373.76"
$ws.Range("E21").Value = "'  -3.24%  "
$ws.Range("D22").Value = "'75.38"
$ws.Range("E22").Value = "'  +1.14%  "
$ws.Range("D23").Value = "'0.562"
$ws.Range("E23").Value = "'  +0.01%  "
$ws.Range("E24").Value = "'  +0.08%  "
$ws.Range("D25").Value = "'3.510.86"
$ws.Range("E25").Value = "'  -2.25%  "
$ws.Range("E26").Value = "'  -5.48%  "
$ws.Range("E27").Value = "'  -3.28%  "
$ws.Range("E28").Value = "'  +0.46%  "
$ws.Range("D29").Value = "'7.38"
$ws.Range("E29").Value = "'  -4.41%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "'  +0.00%  "
$ws.Range("E31").Value = "'  -1.33%  "
$ws.Range("D32").Value = "'7.68"
$ws.Range("E32").Value = "'  -4.89%  "
$ws.Range("D33").Value = "'22.91"
$ws.Range("E33").Value = "'  -1.98%  "
$ws.Range("D34").Value = "'1.30"
$ws.Range("E34").Value = "'  -3.43%  "
$ws.Range("D35").Value = "'5.38"
$ws.Range("E35").Value = "'  +0.31%  "
$ws.Range("D36").Value = "'168.90"
$ws.Range("E36").Value = "'  -0.50%  "
$ws.Range("E37").Value = "'  -4.75%  "
$ws.Range("D38").Value = "'6.77"
$ws.Range("E38").Value = "'  -3.60%  "
$ws.Range("D39").Value = "'29.51"
$ws.Range("E39").Value = "'  -7.51%  "
$ws.Range("D40").Value = "'3.387.39"
$ws.Range("E40").Value = "'  -2.99%  "
$ws.Range("E41").Value = "'  -2.09%  "
$ws.Range("D42").Value = "'42.34"
$ws.Range("E42").Value = "'  -1.23%  "
$ws.Range("D43").Value = "'0.759"
$ws.Range("E43").Value = "'  -4.26%  "
$ws.Range("D44").Value = "'4.31"
$ws.Range("E44").Value = "'  -1.30%  "
$ws.Range("E45").Value = "'  -4.05%  "
$ws.Range("D46").Value = "'1.62"
$ws.Range("E46").Value = "'  -6.04%  "
$ws.Range("D47").Value = "'2.511.19"
$ws.Range("E47").Value = "'  -2.71%  "
$ws.Range("D48").Value = "'22.86"
$ws.Range("E48").Value = "'  +0.99%  "
$ws.Range("E49").Value = "'  -2.47%  "
$ws.Range("E50").Value = "'  +0.02%  "
$ws.Range("D51").Value = "'0.0260"
$ws.Range("E51").Value = "'  -2.50%  "
